$wb = $excel.ActiveWorkbook

# --- Sheet1 ("Sheet1"): remove the consumed name "omp3oay0" from row 2, ---
# --- shifting all the names below it up by one row.                    ---
$names = $wb.Worksheets.Item("Sheet1")
$names.Rows.Item(2).Delete()

# --- Sheet2 ("used"): append a new record for the now-used name,       ---
# --- noting the source file and the timestamp it was used at.          ---
$used = $wb.Worksheets.Item("used")
$nextRow = $used.UsedRange.Rows.Count + 1

$used.Cells.Item($nextRow, 1).Value = "omp3oay0"
$used.Cells.Item($nextRow, 2).Value = "ChatGPT Image 2026年1月24日 01_55_16.png"
$used.Cells.Item($nextRow, 3).Value = "2026-01-24 01:56:28"
